$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the top of the data block
# (row 521), pushing all existing rows down by one. The previously
# last data row (637) now lives at row 638.
$ws.Rows.Item(521).Insert()

$ws.Range("A521").Value = 9
$ws.Range("B521").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C521").Value = "Metropolitana"
$ws.Range("D521").Value = 45211
$ws.Range("E521").Value = 13
$ws.Range("F521").Value = 100112039
$ws.Range("G521").Value = "Ciboulette"
$ws.Range("H521").Value = "Sin especificar"
$ws.Range("I521").Value = "Primera"
$ws.Range("J521").Value = 430
$ws.Range("K521").Value = 1000
$ws.Range("L521").Value = 1100
$ws.Range("M521").Value = 1050
$ws.Range("N521").Value = "`$/docena de atados"
$ws.Range("O521").Value = "Región Metropolitana"
$ws.Range("P521").Value = 350
$ws.Range("Q521").Value = 3
$ws.Range("R521").Value = "Hortaliza"
